$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("contact_info")

# Select the entire row 3 (the blank "   " header placeholder row) and delete it,
# shifting all rows below up by one.
$ws.Rows.Item(3).Select() | Out-Null
$ws.Rows.Item(3).Delete() | Out-Null
